$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Save" header in column H, matching the style of the existing
# header row (bold, bordered, centered) by copying G1's formatting.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New data values for the "Save" column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
